# Update betting odds values on Sheet1 to match the latest FlashScore export.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("Q2").Value = 2.08
$ws.Range("R2").Value = 1.73
$ws.Range("AT2").Value = 2.63

# Row 3
$ws.Range("Q3").Value = 2.35
$ws.Range("R3").Value = 1.57

# Row 5
$ws.Range("Q5").Value = 1.95
$ws.Range("R5").Value = 1.9

# Row 6
$ws.Range("G6").Value = 1.38
$ws.Range("H6").Value = 4.4
$ws.Range("I6").Value = 6.4
$ws.Range("J6").Value = 1.82
$ws.Range("K6").Value = 2.42
$ws.Range("L6").Value = 6
$ws.Range("N6").Value = 13.3
$ws.Range("O6").Value = 1.14
$ws.Range("P6").Value = 4.2
$ws.Range("Q6").Value = 1.55
$ws.Range("R6").Value = 2.15
$ws.Range("S6").Value = 1.28
$ws.Range("T6").Value = 3.34
$ws.Range("U6").Value = 1.85
$ws.Range("V6").Value = 1.91
$ws.Range("W6").Value = 6.9
$ws.Range("X6").Value = 6.2
$ws.Range("Y6").Value = 7.1
$ws.Range("Z6").Value = 7.9
$ws.Range("AA6").Value = 9
$ws.Range("AC6").Value = 13.5
$ws.Range("AD6").Value = 7.8
$ws.Range("AE6").Value = 14.5
$ws.Range("AH6").Value = 16
$ws.Range("AI6").Value = 35
$ws.Range("AJ6").Value = 17
$ws.Range("AN6").Value = 3.3
$ws.Range("AO6").Value = 6.1
$ws.Range("AQ6").Value = 16.5
$ws.Range("AT6").Value = 3.25
$ws.Range("AU6").Value = 7.9
$ws.Range("AW6").Value = 8
$ws.Range("AX6").Value = 37
$ws.Range("AY6").Value = 37

# Row 12
$ws.Range("K12").Value = 1.95

# Row 13
$ws.Range("J13").Value = 1.91
$ws.Range("K13").Value = 2.38
$ws.Range("Q13").Value = 1.88
$ws.Range("R13").Value = 1.98

$wb.Save()
